# Populate header row + append duplicated data block to both sheets.
$wb = $excel.ActiveWorkbook

# ---- CustomerSummary: header row 1 + duplicate rows 2-21 appended as rows 22-41 ----
$ws1 = $wb.Worksheets.Item("CustomerSummary")

$headers1 = @("CustomerID", "FullName", "Email", "DaysSinceRegistration", "LastLoginDate", "IsActive")
for ($j = 0; $j -lt $headers1.Count; $j++) {
    $ws1.Cells.Item(1, $j + 1).Value = $headers1[$j]
}

$sheet1Rows = @(
    @("CUST0001", "FirstName1 LastName1", "customer1@example.com", 476, "2023-10-20 00:00:00", $true),
    @("CUST0002", "FirstName2 LastName2", "customer2@example.com", 725, "2023-02-27 00:00:00", $false),
    @("CUST0003", "FirstName3 LastName3", "customer3@example.com", 481, "2023-10-18 00:00:00", $true),
    @("CUST0004", "FirstName4 LastName4", "customer4@example.com", 761, "2023-01-22 00:00:00", $false),
    @("CUST0005", "FirstName5 LastName5", "customer5@example.com", 595, "2023-07-13 00:00:00", $false),
    @("CUST0006", "FirstName6 LastName6", "customer6@example.com", 458, "2023-11-25 00:00:00", $true),
    @("CUST0007", "FirstName7 LastName7", "customer7@example.com", 581, "2023-07-13 00:00:00", $true),
    @("CUST0008", "FirstName8 LastName8", "customer8@example.com", 733, "2023-02-08 00:00:00", $true),
    @("CUST0009", "FirstName9 LastName9", "customer9@example.com", 489, "2023-10-03 00:00:00", $true),
    @("CUST0010", "FirstName10 LastName10", "customer10@example.com", 602, "2023-07-04 00:00:00", $true),
    @("CUST0011", "FirstName11 LastName11", "customer11@example.com", 703, "2023-03-19 00:00:00", $false),
    @("CUST0012", "FirstName12 LastName12", "customer12@example.com", 702, "2023-03-18 00:00:00", $true),
    @("CUST0013", "FirstName13 LastName13", "customer13@example.com", 668, "2023-04-15 00:00:00", $false),
    @("CUST0014", "FirstName14 LastName14", "customer14@example.com", 594, "2023-07-08 00:00:00", $false),
    @("CUST0015", "FirstName15 LastName15", "customer15@example.com", 409, "2024-01-08 00:00:00", $true),
    @("CUST0016", "FirstName16 LastName16", "customer16@example.com", 448, "2023-12-04 00:00:00", $false),
    @("CUST0017", "FirstName17 LastName17", "customer17@example.com", 714, "2023-03-14 00:00:00", $true),
    @("CUST0018", "FirstName18 LastName18", "customer18@example.com", 698, "2023-03-26 00:00:00", $false),
    @("CUST0019", "FirstName19 LastName19", "customer19@example.com", 727, "2023-02-20 00:00:00", $false),
    @("CUST0020", "FirstName20 LastName20", "customer20@example.com", 425, "2024-01-02 00:00:00", $false)
)

$destStart1 = 22
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $r = $destStart1 + $i
    $row = $sheet1Rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws1.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

# ---- TransactionSummary: header row 1 + duplicate rows 2-20 appended as rows 59-77 ----
$ws2 = $wb.Worksheets.Item("TransactionSummary")

$headers2 = @("CustomerID", "TransactionCount", "TotalAmount", "AverageAmount", "LastTransactionDate", "SuccessRate")
for ($j = 0; $j -lt $headers2.Count; $j++) {
    $ws2.Cells.Item(1, $j + 1).Value = $headers2[$j]
}

$sheet2Rows = @(
    @("CUST0003", 3, 928.3100000000001, 309.4366666666667, "2023-05-15 00:00:00", 0),
    @("CUST0006", 2, 665.53, 332.765, "2023-06-27 00:00:00", 0.5),
    @("CUST0012", 4, 1309.43, 327.3575, "2023-12-04 00:00:00", 0),
    @("CUST0011", 2, 690.4100000000001, 345.205, "2023-11-03 00:00:00", 0),
    @("CUST0018", 4, 2080.15, 520.0375, "2023-10-14 00:00:00", 0.25),
    @("CUST0014", 2, 1303.51, 651.755, "2023-07-06 00:00:00", 0),
    @("CUST0019", 4, 1867.5, 466.875, "2023-09-15 00:00:00", 0.75),
    @("CUST0017", 3, 982.72, 327.5733333333333, "2023-07-04 00:00:00", 0.3333333333333333),
    @("CUST0002", 2, 395.94, 197.97, "2023-05-01 00:00:00", 0.5),
    @("CUST0007", 3, 1438.42, 479.4733333333334, "2023-06-27 00:00:00", 0.3333333333333333),
    @("CUST0001", 6, 3001.88, 500.3133333333333, "2023-12-18 00:00:00", 0.1666666666666667),
    @("CUST0013", 3, 1597.17, 532.39, "2023-10-22 00:00:00", 0),
    @("CUST0005", 2, 1079.01, 539.505, "2023-09-19 00:00:00", 0),
    @("CUST0010", 2, 261.55, 130.775, "2023-06-13 00:00:00", 1),
    @("CUST0004", 3, 820.0700000000001, 273.3566666666667, "2023-11-29 00:00:00", 0.6666666666666666),
    @("CUST0020", 2, 1377.66, 688.83, "2023-07-30 00:00:00", 0),
    @("CUST0008", 1, 330.22, 330.22, "2023-06-07 00:00:00", 0),
    @("CUST0009", 1, 665.04, 665.04, "2023-03-25 00:00:00", 1),
    @("CUST0016", 1, 616.3, 616.3, "2023-04-10 00:00:00", 1)
)

$destStart2 = 59
for ($i = 0; $i -lt $sheet2Rows.Count; $i++) {
    $r = $destStart2 + $i
    $row = $sheet2Rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws2.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
